# Update the "Förändrad" (last-changed) date in column C from 2023-09-05
# (serial 45174) to 2023-09-06 (serial 45175) for every data row in the
# worksheet, as part of the automatic daily refresh of this report.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count + $ws.UsedRange.Row - 1

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45174) {
        $cell.Value2 = 45175
    }
}
